$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 2; existing rows 2-11 shift down to 3-12.
$ws.Rows(2).Insert()

# The insert copies row-1's (header) formatting onto the new row; strip it
# back to the default "no style" formatting used by the other data rows.
$ws.Range("A2:T2").ClearFormats()

# Populate the new row 2 with its data.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C2").Value = "Arica y Parinacota"
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D2").Value = 44972
$ws.Range("E2").Value = 15
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100107
$ws.Range("H2").Value = "Otros"
$ws.Range("I2").Value = 100107011
$ws.Range("J2").Value = "Tuna"
$ws.Range("K2").Value = "Sin especificar"
$ws.Range("L2").Value = "Segunda"
$ws.Range("M2").Value = 140
$ws.Range("N2").Value = 27000
$ws.Range("O2").Value = 28000
$ws.Range("P2").Value = 27429
$ws.Range("Q2").Value = "$/caja 18 kilos"
$ws.Range("R2").Value = "Región Metropolitana"
$ws.Range("S2").Value = 1524
$ws.Range("T2").Value = 18
